$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Save" header in H1 -- copy the formatting used by the other
# header cells (e.g. G1, style index 1: bold/centered/bordered) so the
# new column matches the existing header row, then set its text.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# New "Save" column data values (both rows are 0)
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 0
